# Update summary data and add csv version
#
# Applies the data updates captured in the commit diff to Sheet1 of the
# gftp_summary workbook: new/changed raw measurements in columns C:O, plus
# the corresponding "100/x" ratio formulas in columns T:AF that accompany
# them, and the final active-cell selection left over from editing.
#
# Number formatting for newly-touched ratio cells is picked up from a
# same-style neighbour cell via Copy / PasteSpecial(xlPasteFormats) so the
# shared cellXfs/font/numFmt table entries are reused rather than
# duplicated.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# Reference cells whose formatting we reuse for newly written ratio cells.
$refPurple01 = "T19"   # 0.0 numfmt, purple font  -> style used across rows 19-24
$refPurpleInt = "U19"  # integer numfmt, purple font
$refTheme901 = "U33"   # 0.0 numfmt, theme9 font -> used for style-only tweak
$refTheme902 = "T32"   # 0.00 numfmt, theme9 font -> style used across rows 32-37

# ---------------------------------------------------------------------
# Row 11 block (In Gbps table) -----------------------------------------
$ws.Range("L11").Value = 7.9

# ---------------------------------------------------------------------
# Row 19 (AWS EU C1) ----------------------------------------------------
$ws.Range("N19").Value = 95.2
$ws.Range("O19").Value = 122

$ws.Range($refPurple01).Copy()
$ws.Range("AE19").PasteSpecial($xlPasteFormats)
$ws.Range("AE19").Formula = "=100/N19"

$ws.Range($refPurple01).Copy()
$ws.Range("AF19").PasteSpecial($xlPasteFormats)
$ws.Range("AF19").Formula = "=100/O19"

# ---------------------------------------------------------------------
# Row 20 (Amsterdam) -----------------------------------------------------
$ws.Range("I20").Value = 12.5

$ws.Range($refPurple01).Copy()
$ws.Range("Z20").PasteSpecial($xlPasteFormats)
$ws.Range("Z20").Formula = "=100/I20"

# ---------------------------------------------------------------------
# Row 21 (NewYork) - style of an already-existing formula changes --------
$ws.Range($refPurple01).Copy()
$ws.Range("V21").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 22 (AWS US E1) ------------------------------------------------------
$ws.Range("F22").Value = 3.32

$ws.Range($refPurpleInt).Copy()
$ws.Range("W22").PasteSpecial($xlPasteFormats)

# ---------------------------------------------------------------------
# Row 23 (San Diego) ------------------------------------------------------
$ws.Range("N23").Value = 61.5
$ws.Range("O23").Value = 74.5

$ws.Range($refPurple01).Copy()
$ws.Range("V23").PasteSpecial($xlPasteFormats)
$ws.Range("V23").Formula = "=100/E23"

$ws.Range($refPurple01).Copy()
$ws.Range("AE23").PasteSpecial($xlPasteFormats)
$ws.Range("AE23").Formula = "=100/N23"

$ws.Range($refPurple01).Copy()
$ws.Range("AF23").PasteSpecial($xlPasteFormats)
$ws.Range("AF23").Formula = "=100/O23"

# ---------------------------------------------------------------------
# Row 24 (AWS US W2) -------------------------------------------------------
$ws.Range("L24").Value = 3.3
$ws.Range("N24").Value = 55.4
$ws.Range("O24").Value = 62.5

$ws.Range($refPurpleInt).Copy()
$ws.Range("AC24").PasteSpecial($xlPasteFormats)

$ws.Range($refPurple01).Copy()
$ws.Range("AE24").PasteSpecial($xlPasteFormats)
$ws.Range("AE24").Formula = "=100/N24"

$ws.Range($refPurple01).Copy()
$ws.Range("AF24").PasteSpecial($xlPasteFormats)
$ws.Range("AF24").Formula = "=100/O24"

# ---------------------------------------------------------------------
# Row 32 (AWS EU C1 - In Seconds table) -------------------------------------
$ws.Range("N32").Value = 1810
$ws.Range("O32").Value = 1865

$ws.Range($refTheme902).Copy()
$ws.Range("AE32").PasteSpecial($xlPasteFormats)
$ws.Range("AE32").Formula = "=100/N32"

$ws.Range($refTheme902).Copy()
$ws.Range("AF32").PasteSpecial($xlPasteFormats)
$ws.Range("AF32").Formula = "=100/O32"

# ---------------------------------------------------------------------
# Row 33 (Amsterdam) ----------------------------------------------------
$ws.Range("I33").Value = 194

$ws.Range($refTheme902).Copy()
$ws.Range("Z33").PasteSpecial($xlPasteFormats)
$ws.Range("Z33").Formula = "=100/I33"

# ---------------------------------------------------------------------
# Row 35 (AWS US E1) -----------------------------------------------------
$ws.Range("F35").Value = 16.1

# ---------------------------------------------------------------------
# Row 36 (San Diego) -----------------------------------------------------
$ws.Range("N36").Value = 950
$ws.Range("O36").Value = 995

$ws.Range($refTheme902).Copy()
$ws.Range("AE36").PasteSpecial($xlPasteFormats)
$ws.Range("AE36").Formula = "=100/N36"

$ws.Range($refTheme902).Copy()
$ws.Range("AF36").PasteSpecial($xlPasteFormats)
$ws.Range("AF36").Formula = "=100/O36"

# ---------------------------------------------------------------------
# Row 37 (AWS US W2) -----------------------------------------------------
$ws.Range("F37").Value = 514
$ws.Range("L37").Value = 16.2
$ws.Range("N37").Value = 807
$ws.Range("O37").Value = 900

$ws.Range($refTheme902).Copy()
$ws.Range("W37").PasteSpecial($xlPasteFormats)

$ws.Range($refTheme902).Copy()
$ws.Range("AE37").PasteSpecial($xlPasteFormats)
$ws.Range("AE37").Formula = "=100/N37"

$ws.Range($refTheme902).Copy()
$ws.Range("AF37").PasteSpecial($xlPasteFormats)
$ws.Range("AF37").Formula = "=100/O37"

# ---------------------------------------------------------------------
# Final selection, matching the saved view in the workbook ---------------
$ws.Range("F37").Select()
